$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.113.59"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "'2.476.07"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").Value = "'576.98"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'146.78"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'2.476.33"
$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").Value = "'0.111"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "'29.15"
$ws.Range("E14").Value = "  +8.39%  "

$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "'2.927.35"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").Value = "'63.184.77"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "'2.471.51"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("E22").Value = "  +9.64%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'66.29"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").Value = "'669.24"
$ws.Range("E26").Value = "  +8.58%  "

$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  +14.18%  "

$ws.Range("D28").Value = "'0.0₃0991"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("D31").Value = "'1.45"
$ws.Range("E31").Value = "  +2.96%  "

$ws.Range("D32").Value = "'8.09"
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("E34").Value = "  -2.74%  "

$ws.Range("E35").Value = "  +4.43%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "'4.79"
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'5.45"

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.372"
$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'152.58"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "'18.78"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "'0.0₆0309"
$ws.Range("E45").Value = "  +10.00%  "

$ws.Range("D46").Value = "'150.57"
$ws.Range("E46").Value = "  +4.29%  "

$ws.Range("E47").Value = "  +27.17%  "

$ws.Range("D48").Value = "'3.61"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  +3.22%  "

$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("E51").Value = "  +0.18%  "
